$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.855.54'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '2.085.72'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  +3.11%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.397'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0789'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  +2.79%  '
$ws.Range('D12').Value = '2.392.48'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.79'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').Value = '2.083.34'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '37.789.07'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').Value = '0.0₃0851'
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.66%  '
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0635'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.68'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0985'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.54'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.70%  '
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').Value = '1.451.20'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').Value = '2.277.66'
$ws.Range('E51').Value = '  -0.31%  '
